# 自动更新Excel文件 - 2026-01-30 23:18:20
# Daily countdown update: for each data row, decrement the "剩余" (E, remaining
# days) counter by one. When a row's counter has run out (E = 1), start a new
# cycle: reset E back to the row's total (D) and roll the start date (F)
# forward by 10 days (skip row 36 whose date is malformed and therefore not
# processed by the automation).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    if ($r -eq 36) {
        continue
    }

    $total = $ws.Cells.Item($r, 4).Value2
    $remaining = $ws.Cells.Item($r, 5).Value2
    $start = $ws.Cells.Item($r, 6).Value2

    if ($null -eq $remaining) {
        continue
    }

    if ($remaining -eq 1) {
        $ws.Cells.Item($r, 5).Value = $total
        $ws.Cells.Item($r, 6).Value = $start + 10
    } else {
        $ws.Cells.Item($r, 5).Value = $remaining - 1
    }
}
